$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.105.17'
$ws.Range('E2').Value = '  +2.78%  '

$ws.Range('D3').Value = '2.313.10'
$ws.Range('E3').Value = '  +2.68%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.11'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +5.90%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.535'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.30%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +7.39%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.01'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.44%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.50%  '

$ws.Range('E12').Value = '  +0.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.12'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.32%  '

$ws.Range('D14').Value = '2.672.31'
$ws.Range('E14').Value = '  +2.82%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.98'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.20%  '

$ws.Range('D16').Value = '2.313.35'
$ws.Range('E16').Value = '  +3.20%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.814'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.31%  '

$ws.Range('D18').Value = '43.025.68'
$ws.Range('E18').Value = '  +2.94%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.31%  '

$ws.Range('D20').Value = '0.0₃0920'
$ws.Range('E20').Value = '  +2.31%  '

$ws.Range('E21').Value = '  +3.24%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.40'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.01'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.93%  '

$ws.Range('E25').Value = '  +3.36%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('E27').Value = '  +5.15%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.43'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.66'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.93%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.11'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.18%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.64'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.33'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.95%  '

$ws.Range('E33').Value = '  +0.10%  '

$ws.Range('E34').Value = '  +0.32%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.08'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +7.24%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0743'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.54%  '

$ws.Range('E37').Value = '  +3.12%  '

$ws.Range('E38').Value = '  +0.76%  '

$ws.Range('E39').Value = '  +2.97%  '

$ws.Range('E40').Value = '  +2.16%  '

$ws.Range('E41').Value = '  +7.98%  '

$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.90'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +7.25%  '

$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.31'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.70%  '

$ws.Range('E44').Value = '  +3.48%  '

$ws.Range('D45').Value = '1.974.22'
$ws.Range('E45').Value = '  +0.77%  '

$ws.Range('E46').Value = '  +4.81%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.80'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.75%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.97'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +18.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.76'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.69%  '

$ws.Range('D50').Value = '2.539.22'
$ws.Range('E50').Value = '  +2.67%  '

$ws.Range('E51').Value = '  +4.45%  '
